$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark task "Skapa en liten ikon..." (row 3) as done ("Klar") instead of "Håller på"
$ws.Range("D3").Value = "Klar"

# Sprint day 4 (column I) now has progress on row 3 and in the totals row 6
$ws.Range("I3").Value = 1
$ws.Range("I6").Value = 1

# Update the selected cell as shown in the saved file
$ws.Range("I4").Select()
